# no-op test
$p = $ppt.ActivePresentation
Write-Output $p.Slides.Count
